$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9463879466056824
$ws.Range("B1").Value = 1.634548187255859
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.507455110549927
$ws.Range("E1").Value = 1.339817643165588
